$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "(%) OPERATIVIDAD" column (G) for rows 2-5 switches from numeric 69
# to text-formatted percentage-like codes 99.00 .. 99.03.
$range = $ws.Range("G2:G5")
$range.NumberFormat = "@"

$ws.Range("G2").Value = "99.00"
$ws.Range("G3").Value = "99.01"
$ws.Range("G4").Value = "99.02"
$ws.Range("G5").Value = "99.03"

# Update the active selection to match the author's saved view state.
$ws.Range("G2:G5").Select()
